$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 9
$ws.Range("D9").Value = "파비뉴스 – 정치권의 AI전문가 검증 요청"
$ws.Range("E9").Value = "https://blog.pabii.co.kr/politics-ai-bigdata-specialist/#utm_source=rss&utm_medium=rss&utm_campaign=politics-ai-bigdata-specialist"

# Row 28
$ws.Range("D28").Value = "DQN : Playing Atari with Deep Reinforcement Learning 논문 리뷰 (공부 중)"
$ws.Range("E28").Value = "https://ropiens.tistory.com/75"

# Row 39
$ws.Range("D39").Value = "A Gentle Introduction to Face Recognition in Deep Learning"
$ws.Range("E39").Value = "https://a292run.tistory.com/entry/A-Gentle-Introduction-to-Face-Recognition-in-Deep-Learning-1"
